$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") for all data rows 2..444 from 45190 to 45192
for ($r = 2; $r -le 444; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}

# 2. Row 444 gains an explicit custom row height (matches rest of sheet)
$ws.Rows.Item(444).RowHeight = 15

# 3. Append new row 445: A 44762-2023
$ws.Cells.Item(445, 1).Value = "A 44762-2023"
$ws.Cells.Item(445, 2).Value = 45190
$ws.Cells.Item(445, 3).Value = 45192
$ws.Cells.Item(445, 4).Value = "UPPSALA LÄN"
$ws.Cells.Item(445, 5).Value = "TIERP"
$ws.Cells.Item(445, 7).Value = 0.5
$ws.Cells.Item(445, 8).Value = 0
$ws.Cells.Item(445, 9).Value = 0
$ws.Cells.Item(445, 10).Value = 0
$ws.Cells.Item(445, 11).Value = 0
$ws.Cells.Item(445, 12).Value = 0
$ws.Cells.Item(445, 13).Value = 0
$ws.Cells.Item(445, 14).Value = 0
$ws.Cells.Item(445, 15).Value = 0
$ws.Cells.Item(445, 16).Value = 0
$ws.Cells.Item(445, 17).Value = 0
$ws.Rows.Item(445).RowHeight = 15

# 4. Append new row 446: A 45080-2023
$ws.Cells.Item(446, 1).Value = "A 45080-2023"
$ws.Cells.Item(446, 2).Value = 45191
$ws.Cells.Item(446, 3).Value = 45192
$ws.Cells.Item(446, 4).Value = "UPPSALA LÄN"
$ws.Cells.Item(446, 5).Value = "TIERP"
$ws.Cells.Item(446, 6).Value = "Bergvik skog öst AB"
$ws.Cells.Item(446, 7).Value = 4.5
$ws.Cells.Item(446, 8).Value = 0
$ws.Cells.Item(446, 9).Value = 0
$ws.Cells.Item(446, 10).Value = 0
$ws.Cells.Item(446, 11).Value = 0
$ws.Cells.Item(446, 12).Value = 0
$ws.Cells.Item(446, 13).Value = 0
$ws.Cells.Item(446, 14).Value = 0
$ws.Cells.Item(446, 15).Value = 0
$ws.Cells.Item(446, 16).Value = 0
$ws.Cells.Item(446, 17).Value = 0

# 5. Apply same formatting as the rest of the table:
#    - date format (style 1) for columns B and C
#    - wrap-text empty cell (style 2) for column R
$ws.Range("B445:C446").NumberFormat = $ws.Range("B444:C444").NumberFormat
$ws.Range("R445:R446").WrapText = $true
$ws.Range("R445:R446").Value = ""

# 6. Dimension A1:Y446 is recalculated automatically by the host on save.
